$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at position 58 (pushes the "flip(c[,sn])" row, and everything
# below it, down by one row) and populate it with the new YOCTOPUCE "yset" command.
$ws.Rows.Item(58).Insert()

$ws.Cells.Item(58, 2).Value = "yset(c,b[,sn])"
$ws.Cells.Item(58, 3).Value = "YOCTOPUCE Relay Output: switches channel c of the relay module off (b=0) and on (b=1)"

# Match the row height used by its YOCTOPUCE relay-command neighbours.
$ws.Rows.Item(58).RowHeight = 13.8
